$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
foreach ($addr in @("D4","E4","D5","E5","D6","E6")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D4").Value = "02-03-2026"
$ws.Range("E4").Value = "19-06-2026"
$ws.Range("D5").Value = "08-01-2024"
$ws.Range("E5").Value = "26-04-2024"
$ws.Range("D6").Value = "22-06-2026"
$ws.Range("E6").Value = "22-06-2026"
$ws.Range("F4").Value = 640
$ws.Range("G4").Value = 720
$ws.Range("F5").Value = 80
$ws.Range("G5").Value = 160
$ws.Range("F6").Value = 720
$ws.Range("G6").Value = 720
$ws.Range("I2").Value = "2; 3; 4"
$ws.Range("H4").Value = "1; P3-5"
$ws.Range("H6").Value = "2; 3; 4"

$ws = $wb.Worksheets.Item(2)
foreach ($addr in @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value = "18-09-2023"
$ws.Range("E2").Value = "18-09-2023"
$ws.Range("D3").Value = "29-04-2024"
$ws.Range("E3").Value = "16-08-2024"
$ws.Range("D4").Value = "19-08-2024"
$ws.Range("E4").Value = "06-12-2024"
$ws.Range("D5").Value = "09-12-2024"
$ws.Range("E5").Value = "28-03-2025"
$ws.Range("D6").Value = "31-03-2025"
$ws.Range("E6").Value = "31-03-2025"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("F3").Value = 160
$ws.Range("G3").Value = 240
$ws.Range("F4").Value = 240
$ws.Range("G4").Value = 320
$ws.Range("F5").Value = 320
$ws.Range("G5").Value = 400
$ws.Range("F6").Value = 400
$ws.Range("G6").Value = 400
$ws.Range("I2").Value = "2; 3; 4"
$ws.Range("H6").Value = "2; 3; 4"

$ws = $wb.Worksheets.Item(3)
foreach ($addr in @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value = "18-09-2023"
$ws.Range("E2").Value = "18-09-2023"
$ws.Range("D3").Value = "31-03-2025"
$ws.Range("E3").Value = "18-07-2025"
$ws.Range("D4").Value = "21-07-2025"
$ws.Range("E4").Value = "07-11-2025"
$ws.Range("D5").Value = "10-11-2025"
$ws.Range("E5").Value = "27-02-2026"
$ws.Range("D6").Value = "02-03-2026"
$ws.Range("E6").Value = "02-03-2026"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("F3").Value = 400
$ws.Range("G3").Value = 480
$ws.Range("F4").Value = 480
$ws.Range("G4").Value = 560
$ws.Range("F5").Value = 560
$ws.Range("G5").Value = 640
$ws.Range("F6").Value = 640
$ws.Range("G6").Value = 640
$ws.Range("I2").Value = "2; 3; 4"
$ws.Range("H6").Value = "2; 3; 4"

$ws = $wb.Worksheets.Item(4)
foreach ($addr in @("D4","E4","D6","E6")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D4").Value = "08-01-2024"
$ws.Range("E4").Value = "26-04-2024"
$ws.Range("D6").Value = "29-04-2024"
$ws.Range("E6").Value = "29-04-2024"
$ws.Range("F4").Value = 80
$ws.Range("G4").Value = 160
$ws.Range("F6").Value = 160
$ws.Range("G6").Value = 160
$ws.Range("I2").Value = "2; 3; 4"
$ws.Range("H4").Value = "1; P3-5"
$ws.Range("H6").Value = "2; 3; 4"

$ws = $wb.Worksheets.Item(5)
foreach ($addr in @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value = "18-09-2023"
$ws.Range("E2").Value = "18-09-2023"
$ws.Range("D3").Value = "18-09-2023"
$ws.Range("E3").Value = "05-01-2024"
$ws.Range("D4").Value = "18-09-2023"
$ws.Range("E4").Value = "05-01-2024"
$ws.Range("D5").Value = "18-09-2023"
$ws.Range("E5").Value = "05-01-2024"
$ws.Range("D6").Value = "08-01-2024"
$ws.Range("E6").Value = "08-01-2024"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 80
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 80
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 80
$ws.Range("F6").Value = 80
$ws.Range("G6").Value = 80
$ws.Range("I2").Value = "2; 3; 4"
$ws.Range("H6").Value = "2; 3; 4"

$ws = $wb.Worksheets.Item(6)
foreach ($addr in @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6")) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value = "18-09-2023"
$ws.Range("E2").Value = "18-09-2023"
$ws.Range("D3").Value = "18-09-2023"
$ws.Range("E3").Value = "05-01-2024"
$ws.Range("D4").Value = "18-09-2023"
$ws.Range("E4").Value = "05-01-2024"
$ws.Range("D5").Value = "18-09-2023"
$ws.Range("E5").Value = "05-01-2024"
$ws.Range("D6").Value = "08-01-2024"
$ws.Range("E6").Value = "08-01-2024"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 80
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 80
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 80
$ws.Range("F6").Value = 80
$ws.Range("G6").Value = 80
$ws.Range("I2").Value = "2; 3; 4"
$ws.Range("H6").Value = "2; 3; 4"
